$d = $word.ActiveDocument

# 1. Update the date paragraph.
$d.Content.Find.Execute("2024-04-05", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-07", 2)

# 2. Remove the pie chart figure, its source link, the four aggregated bar
#    plot figures, and their trailing source link -- i.e. everything after
#    the abstract paragraph through the end of the body (before the
#    sectPr).
$startPara = $d.Paragraphs.Item(7)
$endPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)
$rng.Delete()
